# Automatic Excel update [2025-07-29 07:28:03]
# Applies the daily-refresh edit to wyniki_ofert_z_filtra.xlsx:
#  - bumps "Data ostatniej aktualizacji" (col E) from 2025-07-28 -> 2025-07-29
#    for every still-active listing row
#  - flips "Aktywne" (col H) from TRUE to FALSE for listings that were
#    replaced/relisted (their col E date is NOT bumped, matching the
#    already-archived rows at the bottom of each sheet)
#  - appends newly found listings as new rows at the bottom of each sheet

$wb = $excel.ActiveWorkbook

function Set-TextValue($cell, [string]$text) {
    # Force literal text so Excel doesn't auto-coerce yyyy-mm-dd-looking
    # strings into a date serial, then restore the default "Normal" style
    # so the cell keeps looking like every other untouched text cell
    # (no lingering custom number format on the cell itself).
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------
# Sheet 1: "powiat krakowski"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("powiat krakowski")

# Rows whose listing got delisted/replaced -> Aktywne goes to FALSE and
# their "last updated" date is left untouched (same pattern as the
# archived "/hpr/" rows already in the sheet).
$ws1DeactivatedRows = @(12, 45)
foreach ($r in $ws1DeactivatedRows) {
    $ws1.Cells.Item($r, 8).Value = $false
}

# All remaining data rows (2..59) whose update date is still the old
# "2025-07-28" get bumped to "2025-07-29". (A couple of rows already sit
# on an older/stale date and are left alone, matching upstream.)
for ($r = 2; $r -le 59; $r++) {
    if ($ws1DeactivatedRows -contains $r) { continue }
    $cell = $ws1.Cells.Item($r, 5)
    if ($cell.Text -eq "2025-07-28") {
        Set-TextValue $cell "2025-07-29"
    }
}

# New listings found today.
$ws1.Cells.Item(61, 1).Value = "Działka blisko Krakowa w pięknej okolicy!"
$ws1.Cells.Item(61, 2).Value = "Kamień, Czernichów, krakowski, małopolskie"
$ws1.Cells.Item(61, 3).Value = 248000
Set-TextValue $ws1.Cells.Item(61, 4) "2025-07-29"
Set-TextValue $ws1.Cells.Item(61, 5) "2025-07-29"
$ws1.Cells.Item(61, 6).Value = 248000
$ws1.Cells.Item(61, 7).Value = 0.49
$ws1.Cells.Item(61, 8).Value = $true
$ws1.Cells.Item(61, 9).Value = "https://www.otodom.pl/pl/oferta/dzialka-blisko-krakowa-w-pieknej-okolicy-ID4xcWR"

$ws1.Cells.Item(62, 1).Value = "Mogilany - działka dla miłośników ogrodów. Polecam"
$ws1.Cells.Item(62, 2).Value = "ul. Górska, Mogilany, Mogilany, krakowski, małopolskie"
$ws1.Cells.Item(62, 3).Value = 250000
Set-TextValue $ws1.Cells.Item(62, 4) "2025-07-29"
Set-TextValue $ws1.Cells.Item(62, 5) "2025-07-29"
$ws1.Cells.Item(62, 6).Value = 250000
$ws1.Cells.Item(62, 7).Value = 0.49
$ws1.Cells.Item(62, 8).Value = $true
$ws1.Cells.Item(62, 9).Value = "https://www.otodom.pl/pl/oferta/mogilany-dzialka-dla-milosnikow-ogrodow-polecam-ID4w22e"

# ---------------------------------------------------------------------
# Sheet 2: "powiat wielicki"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("powiat wielicki")

$ws2DeactivatedRows = @(17, 40, 54)
foreach ($r in $ws2DeactivatedRows) {
    $ws2.Cells.Item($r, 8).Value = $false
}

for ($r = 2; $r -le 55; $r++) {
    if ($ws2DeactivatedRows -contains $r) { continue }
    $cell = $ws2.Cells.Item($r, 5)
    if ($cell.Text -eq "2025-07-28") {
        Set-TextValue $cell "2025-07-29"
    }
}

# New listings found today.
$ws2.Cells.Item(56, 1).Value = "Działka budowlana 20a, koło Wieliczki"
$ws2.Cells.Item(56, 2).Value = "Dobranowice, Wieliczka, wielicki, małopolskie"
$ws2.Cells.Item(56, 3).Value = 245000
Set-TextValue $ws2.Cells.Item(56, 4) "2025-07-29"
Set-TextValue $ws2.Cells.Item(56, 5) "2025-07-29"
$ws2.Cells.Item(56, 6).Value = 245000
$ws2.Cells.Item(56, 7).Value = 0.49
$ws2.Cells.Item(56, 8).Value = $true
$ws2.Cells.Item(56, 9).Value = "https://www.otodom.pl/pl/oferta/dzialka-budowlana-20a-kolo-wieliczki-ID4xseG"

$ws2.Cells.Item(57, 1).Value = "Działka budowlana 20a, koło Wieliczki"
$ws2.Cells.Item(57, 2).Value = "Dobranowice, Wieliczka, wielicki, małopolskie"
$ws2.Cells.Item(57, 3).Value = 245000
Set-TextValue $ws2.Cells.Item(57, 4) "2025-07-29"
Set-TextValue $ws2.Cells.Item(57, 5) "2025-07-29"
$ws2.Cells.Item(57, 6).Value = 245000
$ws2.Cells.Item(57, 7).Value = 0.49
$ws2.Cells.Item(57, 8).Value = $true
$ws2.Cells.Item(57, 9).Value = "https://www.otodom.pl/hpr/pl/oferta/dzialka-budowlana-20a-kolo-wieliczki-ID4xseG"
